$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Class Statistics block (K/L columns, rows 6-10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 141

$ws.Range("G8").Value = "dnasr281@gmail.com, System"
$ws.Range("L8").Value = 168

$ws.Range("G9").Value = "dnasr281@gmail.com, System"
$ws.Range("L9").Value = "44.3%"

$ws.Range("G10").Value = "dnasr281@gmail.com, System"
$ws.Range("L10").Value = "72.6%"

# ---------------------------------------------------------------------------
# Row 14 (B1A1 / session 14) flips from "Pending" (yellow) to "Recorded"
# (green) - reuse the formatting already applied to a "Recorded" row (row 6)
# ---------------------------------------------------------------------------
$ws.Range("A6:I6").Copy()
$ws.Range("A14:I14").PasteSpecial(-4122)
$ws.Range("G14").Value = "dnasr281@gmail.com"
$ws.Range("H14").Value = "24/26"
$ws.Range("I14").Value = "Recorded"

# Group Statistics table (rows 15-20) reflecting the newly recorded session
$ws.Range("O15").Value = 11
$ws.Range("Q15").Value = 13
$ws.Range("R15").Value = "42.3%"
$ws.Range("S15").Value = "81.8%"

$ws.Range("O16").Value = 12
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = "46.2%"
$ws.Range("S16").Value = "78.4%"

$ws.Range("O17").Value = 12
$ws.Range("Q17").Value = 13
$ws.Range("R17").Value = "46.2%"
$ws.Range("S17").Value = "64.7%"

$ws.Range("O18").Value = 12
$ws.Range("Q18").Value = 13
$ws.Range("R18").Value = "46.2%"
$ws.Range("S18").Value = "68.5%"

$ws.Range("O19").Value = 12
$ws.Range("Q19").Value = 13
$ws.Range("R19").Value = "46.2%"
$ws.Range("S19").Value = "70.8%"

$ws.Range("O20").Value = 11
$ws.Range("Q20").Value = 13
$ws.Range("R20").Value = "42.3%"
$ws.Range("S20").Value = "74.3%"

# ---------------------------------------------------------------------------
# B1A2 section
# ---------------------------------------------------------------------------
$ws.Range("G34").Value = "dnasr281@gmail.com, System"
$ws.Range("G35").Value = "dnasr281@gmail.com, System"
$ws.Range("G36").Value = "dnasr281@gmail.com, System"

$ws.Range("A6:I6").Copy()
$ws.Range("A40:I40").PasteSpecial(-4122)
$ws.Range("G40").Value = "dnasr281@gmail.com"
$ws.Range("H40").Value = "24/27"
$ws.Range("I40").Value = "Recorded"

# ---------------------------------------------------------------------------
# B1B1 section
# ---------------------------------------------------------------------------
$ws.Range("G60").Value = "dnasr281@gmail.com, System"
$ws.Range("G61").Value = "dnasr281@gmail.com, System"
$ws.Range("G62").Value = "dnasr281@gmail.com, System"

$ws.Range("A6:I6").Copy()
$ws.Range("A66:I66").PasteSpecial(-4122)
$ws.Range("G66").Value = "dnasr281@gmail.com"
$ws.Range("H66").Value = "20/26"
$ws.Range("I66").Value = "Recorded"

# ---------------------------------------------------------------------------
# B1B2 section
# ---------------------------------------------------------------------------
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("G87").Value = "dnasr281@gmail.com, System"
$ws.Range("G88").Value = "dnasr281@gmail.com, System"

$ws.Range("A6:I6").Copy()
$ws.Range("A92:I92").PasteSpecial(-4122)
$ws.Range("G92").Value = "dnasr281@gmail.com"
$ws.Range("H92").Value = "23/27"
$ws.Range("I92").Value = "Recorded"

# ---------------------------------------------------------------------------
# B1C1 section
# ---------------------------------------------------------------------------
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G113").Value = "dnasr281@gmail.com, System"
$ws.Range("G114").Value = "dnasr281@gmail.com, System"

$ws.Range("A6:I6").Copy()
$ws.Range("A118:I118").PasteSpecial(-4122)
$ws.Range("G118").Value = "dnasr281@gmail.com"
$ws.Range("H118").Value = "26/30"
$ws.Range("I118").Value = "Recorded"

# ---------------------------------------------------------------------------
# B1C2 section
# ---------------------------------------------------------------------------
$ws.Range("G138").Value = "dnasr281@gmail.com, System"
$ws.Range("G139").Value = "dnasr281@gmail.com, System"
$ws.Range("G140").Value = "dnasr281@gmail.com, System"

$ws.Range("A6:I6").Copy()
$ws.Range("A144:I144").PasteSpecial(-4122)
$ws.Range("G144").Value = "dnasr281@gmail.com"
$ws.Range("H144").Value = "20/23"
$ws.Range("I144").Value = "Recorded"

# ---------------------------------------------------------------------------
# Remaining simple "System, X" -> "X, System" swaps (B1D1, B1D2, B1E1,
# B1E2, B1F1, B1F2 sections)
# ---------------------------------------------------------------------------
$ws.Range("G164").Value = "dnasr281@gmail.com, System"
$ws.Range("G167").Value = "dnasr281@gmail.com, System"

$ws.Range("G191").Value = "dnasr281@gmail.com, System"
$ws.Range("G194").Value = "dnasr281@gmail.com, System"

$ws.Range("G218").Value = "dnasr281@gmail.com, System"
$ws.Range("G221").Value = "dnasr281@gmail.com, System"

$ws.Range("G245").Value = "dnasr281@gmail.com, System"
$ws.Range("G248").Value = "dnasr281@gmail.com, System"

$ws.Range("G272").Value = "dnasr281@gmail.com, System"
$ws.Range("G275").Value = "dnasr281@gmail.com, System"

$ws.Range("G299").Value = "dnasr281@gmail.com, System"
$ws.Range("G302").Value = "dnasr281@gmail.com, System"
